$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes existing rows 5..39 down to 6..40)
# for the newly reported day, matching the source table's pattern of
# prepending each day's entry under the totals row.
$ws.Rows(5).Insert()

# Match the formatting of the data rows (use the row right below, which is
# the row that used to be row 5, as the formatting template).
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)

# New row 5: data for 2021-05-24 (Mon)
$ws.Range("A5").Value = 44340
$ws.Range("B5").Value = "(月)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 136760
$ws.Range("E5").Value = 193829

# Update the running-total row (row 4) to include the new day's figures
$ws.Range("D4").Value = 4102171
$ws.Range("E4").Value = 2666805

# Update the "as of" date label
$ws.Range("E2").Value = "（5月24日時点）"
